# Aula 16 - Habitos de Performance
# Fix: the logo image lived only on the shared slide layout (slideLayout2 /
# "MASTER"), which some PowerPoint builds fail to open. The fix embeds the
# logo picture + the dark background fill directly on every slide, and
# re-points every slide at the plain "DEFAULT" layout so the now-unused
# "MASTER" layout can be dropped from the deck.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Text tweaks (shorten a handful of bullet / title strings).
#    Shape indices below are the *pre-picture-insertion* indices (the order
#    the shapes already sit in on each slide).
# ---------------------------------------------------------------------------

# Slide 2 - "Objetivos da Aula"
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(4).TextFrame.TextRange.Text = "Entender ciência dos hábitos"
$s2.Shapes.Item(6).TextFrame.TextRange.Text = "Criar novos hábitos"

# Slide 3 - "Ciência dos Hábitos"
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "Ciência"
$s3.Shapes.Item(6).TextFrame.TextRange.Text = "Economizam energia"

# Slide 4 - "Criando Novos Hábitos"
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "Criar Novos"
$s4.Shapes.Item(4).TextFrame.TextRange.Text = "Comece pequeno"
$s4.Shapes.Item(6).TextFrame.TextRange.Text = "Encadeamento"
$s4.Shapes.Item(10).TextFrame.TextRange.Text = "Recompensas"

# Slide 5 - "Quebrando Hábitos Ruins"
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "Quebrar Ruins"

# Slide 6 - "Acompanhamento"
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(8).TextFrame.TextRange.Text = "Aceite falhas"

# Slide 7 - "Pontos-Chave"
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(4).TextFrame.TextRange.Text = "Hábitos > motivação"
$s7.Shapes.Item(8).TextFrame.TextRange.Text = "Comece pequeno, encadeie"
$s7.Shapes.Item(10).TextFrame.TextRange.Text = "66 dias para formar"

# Slide 8 - "Atividade Prática"
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(4).TextFrame.TextRange.Text = "Implantar 1 Hábito: loop + tracker + 2 semanas"

# Slide 9 - "Aula 16 - Encerramento"
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(4).TextFrame.TextRange.Text = """Você cai ao nível dos seus sistemas."""

# ---------------------------------------------------------------------------
# 2. Give every slide its own copy of the logo + dark background, then point
#    it at the plain DEFAULT layout instead of the logo-carrying MASTER one.
# ---------------------------------------------------------------------------

$master = $p.SlideMaster
$defaultLayout = $master.CustomLayouts.Item(1)   # "DEFAULT" (no logo)
$masterLayout  = $master.CustomLayouts.Item(2)   # "MASTER" (logo only lives here)

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    # Dark navy background, set directly on the slide.
    $slide.FollowMasterBackground = $false
    $slide.Background.Fill.Solid()
    $slide.Background.Fill.ForeColor.RGB = 0x2E1A1A

    # Embed the logo image on the slide itself (points = EMU / 12700).
    $logo = $slide.Shapes.AddPicture("preencoded.png", $false, $true, 21.6, 10.8, 86.4, 36.0)
    $logo.Name = "Image 0"
    $logo.AlternativeText = "preencoded.png"
    $logo.LockAspectRatio = -1
    $logo.ZOrder(1)

    # Slide no longer needs the logo-only layout.
    $slide.CustomLayout = $defaultLayout
}

# Nothing references the "MASTER" layout anymore - drop it from the deck.
$masterLayout.Delete()
